$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new values look like plain decimals (e.g. "36.30").
# Pin those specific cells to text format first so Excel keeps the exact
# string (incl. trailing zeros) instead of silently coercing to a Number.
$textCells = @("D10", "D11", "D12", "D14", "D17", "D19", "D20", "D23", "D24", "D26", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D47", "D49", "D5", "D50", "D51", "D6", "D9")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '49.891.70'
$ws.Range("E2").Value = '  +4.19%  '
$ws.Range("D3").Value = '2.675.03'
$ws.Range("E3").Value = '  +7.70%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '113.56'
$ws.Range("E5").Value = '  +8.73%  '
$ws.Range("D6").Value = '326.15'
$ws.Range("E6").Value = '  +2.83%  '
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  +3.35%  '
$ws.Range("D10").Value = '40.77'
$ws.Range("E10").Value = '  +5.00%  '
$ws.Range("D11").Value = '20.12'
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").Value = '0.0823'
$ws.Range("E12").Value = '  +3.22%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '7.36'
$ws.Range("E14").Value = '  +4.85%  '
$ws.Range("D15").Value = '3.096.78'
$ws.Range("D16").Value = '2.676.27'
$ws.Range("E16").Value = '  +9.14%  '
$ws.Range("D17").Value = '0.876'
$ws.Range("E17").Value = '  +6.46%  '
$ws.Range("D18").Value = '49.856.41'
$ws.Range("E18").Value = '  +4.20%  '
$ws.Range("D19").Value = '13.17'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").Value = '6.78'
$ws.Range("E20").Value = '  +4.23%  '
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").Value = '0.0₃0960'
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '278.77'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '71.91'
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").Value = '26.84'
$ws.Range("E26").Value = '  +4.84%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +5.98%  '
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").Value = '36.30'
$ws.Range("E30").Value = '  +5.36%  '
$ws.Range("D31").Value = '0.142'
$ws.Range("E31").Value = '  +3.26%  '
$ws.Range("D32").Value = '50.30'
$ws.Range("E32").Value = '  +2.41%  '
$ws.Range("D33").Value = '5.47'
$ws.Range("E33").Value = '  +4.39%  '
$ws.Range("D34").Value = '19.49'
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("D35").Value = '0.0809'
$ws.Range("E35").Value = '  +5.32%  '
$ws.Range("D36").Value = '5.10'
$ws.Range("E36").Value = '  +13.75%  '
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '2.07'
$ws.Range("E38").Value = '  +7.30%  '
$ws.Range("D39").Value = '3.15'
$ws.Range("E39").Value = '  +10.23%  '
$ws.Range("D40").Value = '125.29'
$ws.Range("E40").Value = '  +2.10%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.113'
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '22.59'
$ws.Range("E42").Value = '  +5.58%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").Value = '  +6.46%  '
$ws.Range("D45").Value = '2.126.09'
$ws.Range("E45").Value = '  +7.09%  '
$ws.Range("E46").Value = '  +6.64%  '
$ws.Range("D47").Value = '2.24'
$ws.Range("E47").Value = '  +11.75%  '
$ws.Range("E48").Value = '  +8.01%  '
$ws.Range("D49").Value = '9.04'
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("D50").Value = '5.35'
$ws.Range("E50").Value = '  +5.67%  '
$ws.Range("D51").Value = '59.38'
$ws.Range("E51").Value = '  +6.83%  '
